# Generate Report for Handoff
# - Update status text from "In Translation" to "Ready for handoff"
# - Bump the handoff timestamps by 30 seconds
# - Re-fit the affected columns so the new, longer status text is not truncated

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Status column(s): "In Translation" -> "Ready for handoff"
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# Handoff datetimes move forward by 30 seconds
$overview.Range("G2").Value = "2016-09-06 23:12:57"
$dede.Range("H2").Value = "2016-09-06 23:12:57"
$zhcn.Range("H2").Value = "2016-09-06 23:12:52"

# Re-fit the columns that now hold the longer "Ready for handoff" text
$overview.Columns.Item(5).ColumnWidth = 17.2159881591797
$overview.Columns.Item(6).ColumnWidth = 17.2159881591797
$zhcn.Columns.Item(3).ColumnWidth = 17.2159881591797
$dede.Columns.Item(3).ColumnWidth = 17.2159881591797
